# Better Pressed Buttons (UI Enhancement)
# Update the freelancer schedule values so the "pressed" (active) shift
# button per person/day matches the corrected schedule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B3" = "off";   "C3" = "15-24"; "D3" = "off";   "E3" = "off";  "F3" = "7-16";  "G3" = "10-19"
    "B4" = "off";   "C4" = "15-24"; "D4" = "15-24"; "E4" = "10-19"; "F4" = "off";  "G4" = "7-16"
    "B5" = "15-24"; "C5" = "off";   "D5" = "off";   "E5" = "15-24"; "F5" = "10-19"; "G5" = "7-16"
    "B6" = "15-24"; "C6" = "15-24"; "D6" = "off";   "E6" = "off";  "F6" = "off";   "G6" = "7-16"
    "B7" = "7-16";  "C7" = "10-19"; "D7" = "15-24"; "E7" = "off";  "F7" = "off";   "G7" = "15-24"
    "B8" = "7-16";  "C8" = "off";   "D8" = "10-19"; "E8" = "off";  "F8" = "15-24"; "G8" = "15-24"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
